$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Catálogo de Mudança" ---
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 350 (pushes old rows 350..456 down to 351..457)
$ws.Rows.Item(350).Insert()

# Copy the formatting from row 349 (an existing "Alta" priority row) onto
# the freshly inserted row 350 so styles (borders/fills/fonts) match.
$ws.Range("A349:F349").Copy()
$ws.Range("A350:F350").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new box entry: "Caixa 42 / Acessórios / Bijuterias / Suíte / Alta"
$ws.Range("A350").Value = "Caixa 42"
$ws.Range("B350").Value = "Acessórios"
$ws.Range("C350").Value = "Bijuterias"
$ws.Range("D350").Value = "Suíte"
$ws.Range("E350").Value = "Alta"

# --- Sheet 2: "Resumo" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B3").Value = 456   # Total de Itens
$ws2.Range("B6").Value = 105   # Alta
